$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1.480"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "12.579"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "5.547 %"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "5.624 %"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "5.848 %"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "7.099"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "34.332"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "10.464"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "4.614 %"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "10.137"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.509"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "1.106 %"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "13.062"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "24.016"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "10.591 %"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "13.172"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "15.122"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "9.952"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "4.389 %"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "15.261"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "18.188"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "14.765"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "6.511 %"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "18.188"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "10.056"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "4.435 %"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "16.755"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "7.389 %"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "25.017"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "10.345"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "4.562 %"
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "25.456"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "25.018"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.526"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "13.218"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "5.829 %"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "18.811"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "8.295 %"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "15.613"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "6.885 %"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "40.301"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.976"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15.123"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "6.669 %"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "42.314"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "41.389"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "44.928"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "40.537"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "26.546"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "11.706 %"
